$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row at row 11 (shifts existing rows 11+ down by one)
$ws.Rows("11:11").Insert()

# Populate the new row with the "Nickname" column-index setting
$ws.Range("A11").Value = "BirthdayList_Index_Nickname"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "Column index of ""Nickname"" in BirthdayList datatable"

# Update the view state to match the recorded selection after the edit
$null = $ws.Range("C12").Select()
